$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 44523.82436342593
$ws.Range("B4").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("C4").Value = 44523.82605324074
$ws.Range("C4").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("D4").Value = "IP Address"
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 145
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 44523.82605324074
$ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("I4").Value = "3mihar"
$ws.Range("J4").Value = "ebola %>% `npivot_longer(cols = Cases_Guinea:last_col()) %>% `nseparate(name, into = c(`"case_death`", `"country`")) %>% `ndrop_na()"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44523.82412037037
$ws.Range("B5").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("C5").Value = 44523.82731481482
$ws.Range("C5").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("D5").Value = "Spam"
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 275
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 44523.82731481482
$ws.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("I5").Value = "1pogus"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44523.09984953704
$ws.Range("B6").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("C6").Value = 44523.10025462963
$ws.Range("C6").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("D6").Value = "Spam"
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 34
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 44523.87974537037
$ws.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("I6").Value = "2nesch"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44518.84451388889
$ws.Range("B7").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("C7").Value = 44518.85037037038
$ws.Range("C7").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("D7").Value = "IP Address"
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 505
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 44523.8797800926
$ws.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("I7").Value = "0hacar"

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 44518.82807870371
$ws.Range("B8").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("C8").Value = 44518.82834490741
$ws.Range("C8").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("D8").Value = "Spam"
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 22
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 44523.87981481482
$ws.Range("H8").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("I8").Value = "3kusou"
